$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64
$ws.Range("H64").Value = 9057.058999999999
$ws.Range("I64").Value = 3424.4285
$ws.Range("J64").Value = 12999.9
$ws.Range("K64").Value = 3424.4285
$ws.Range("L64").Value = 12999.9
$ws.Range("M64").Value = -3176.4285
$ws.Range("N64").Value = -13495.9
# Row 67
$ws.Range("H67").Value = 9057.058999999999
$ws.Range("I67").Value = 3424.4285
$ws.Range("J67").Value = 12999.9
$ws.Range("K67").Value = 3424.4285
$ws.Range("L67").Value = 12999.9
$ws.Range("M67").Value = -2566.4285
$ws.Range("N67").Value = -14715.9
# Row 98
$ws.Range("H98").Value = 939.5
$ws.Range("J98").Value = 636
$ws.Range("L98").Value = 636
$ws.Range("N98").Value = -3632
# Row 122
$ws.Range("H122").Value = 939.5
$ws.Range("J122").Value = 636
$ws.Range("L122").Value = 1908
$ws.Range("N122").Value = -6808
# Row 135
$ws.Range("H135").Value = 534
$ws.Range("I135").Value = 445.16666
$ws.Range("K135").Value = 4006.49994
$ws.Range("M135").Value = -1471.49994
# Row 138
$ws.Range("H138").Value = 1984.2545
$ws.Range("J138").Value = 2366.5833
$ws.Range("L138").Value = 7099.749899999999
$ws.Range("N138").Value = -17379.7499

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2064
$ws.Range("J2").Value = 4722.857
$ws.Range("L2").Value = 4722.857
$ws.Range("N2").Value = -4948.857
# Row 74
$ws.Range("H74").Value = 3047.4666
$ws.Range("J74").Value = 4965.3076
$ws.Range("L74").Value = 4965.3076
$ws.Range("N74").Value = -6713.3076
# Row 77
$ws.Range("H77").Value = 3047.4666
$ws.Range("J77").Value = 4965.3076
$ws.Range("L77").Value = 24826.538
$ws.Range("N77").Value = -33562.538
# Row 116
$ws.Range("H116").Value = 2064
$ws.Range("J116").Value = 4722.857
$ws.Range("L116").Value = 4722.857
$ws.Range("N116").Value = -9310.857
# Row 122
$ws.Range("H122").Value = 3420.7058
$ws.Range("I122").Value = 3648.7778
$ws.Range("J122").Value = 3164.125
$ws.Range("K122").Value = 10946.3334
$ws.Range("L122").Value = 9492.375
$ws.Range("M122").Value = -8496.3334
$ws.Range("N122").Value = -14392.375
# Row 132
$ws.Range("H132").Value = 4230.171
$ws.Range("I132").Value = 3665.2
$ws.Range("K132").Value = 10995.6
$ws.Range("M132").Value = -8465.599999999999

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2064
$ws.Range("J3").Value = 4722.857
$ws.Range("L3").Value = 4722.857
$ws.Range("N3").Value = -4950.857
# Row 86
$ws.Range("H86").Value = 1333.2222
$ws.Range("I86").Value = 1302.7333
$ws.Range("K86").Value = 1302.7333
$ws.Range("M86").Value = -179.7333000000001
# Row 89
$ws.Range("H89").Value = 1333.2222
$ws.Range("I89").Value = 1302.7333
$ws.Range("K89").Value = 6513.6665
$ws.Range("M89").Value = -897.6665000000003
# Row 128
$ws.Range("H128").Value = 2500
$ws.Range("I128").Value = 2500
$ws.Range("K128").Value = 7500
$ws.Range("M128").Value = -5010

$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 3608.111
$ws.Range("I122").Value = 3598.2856
$ws.Range("K122").Value = 10794.8568
$ws.Range("M122").Value = -8344.856800000001

$ws = $wb.Worksheets.Item("CUL")
# Row 55
$ws.Range("H55").Value = 824.6
$ws.Range("J55").Value = 879.75
$ws.Range("L55").Value = 2639.25
$ws.Range("N55").Value = -2993.25
# Row 117
$ws.Range("H117").Value = 2502605.5
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 2502605.5
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 7507816.5
$ws.Range("M117").Value = ""
$ws.Range("N117").Value = -7514700.5
# Row 123
$ws.Range("H123").Value = 13899.429
$ws.Range("I123").Value = 4824.75
$ws.Range("J123").Value = 25999
$ws.Range("K123").Value = 14474.25
$ws.Range("L123").Value = 77997
$ws.Range("M123").Value = -12024.25
$ws.Range("N123").Value = -82897
# Row 128
$ws.Range("H128").Value = 691178.9399999999
$ws.Range("I128").Value = 691178.9399999999
$ws.Range("K128").Value = 2073536.82
$ws.Range("M128").Value = -2068556.82
# Row 131
$ws.Range("H131").Value = 21742596
$ws.Range("I131").Value = 166667660
$ws.Range("J131").Value = 3834.6
$ws.Range("K131").Value = 500002980
$ws.Range("L131").Value = 11503.8
$ws.Range("M131").Value = -499997940
$ws.Range("N131").Value = -21583.8

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 10140.8
$ws.Range("I2").Value = 156.55556
$ws.Range("J2").Value = 99999
$ws.Range("K2").Value = 156.55556
$ws.Range("L2").Value = 99999
$ws.Range("M2").Value = -43.55556000000001
$ws.Range("N2").Value = -100225
# Row 18
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = ""
# Row 19
$ws.Range("H19").Value = 13601.333
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 13601.333
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 13601.333
$ws.Range("M19").Value = -14177.333
# Row 92
$ws.Range("H92").Value = 1650
$ws.Range("J92").Value = 1650
$ws.Range("L92").Value = 1650
$ws.Range("N92").Value = -5394
# Row 113
$ws.Range("H113").Value = 1094.5714
$ws.Range("I113").Value = 1072.5
$ws.Range("J113").Value = 1124
$ws.Range("K113").Value = 1072.5
$ws.Range("L113").Value = 1124
$ws.Range("M113").Value = 1097.5
$ws.Range("N113").Value = -5464
# Row 132
$ws.Range("H132").Value = 2487.7334
$ws.Range("I132").Value = 2451.1428
$ws.Range("K132").Value = 7353.428400000001
$ws.Range("M132").Value = -4823.428400000001

$ws = $wb.Worksheets.Item("LTW")
# Row 58
$ws.Range("I58").Value = 9331.666999999999
$ws.Range("J58").Value = 20000
$ws.Range("K58").Value = 9331.666999999999
$ws.Range("L58").Value = 20000
$ws.Range("M58").Value = -9071.666999999999
$ws.Range("N58").Value = -20520
# Row 101
$ws.Range("H101").Value = 64489.5
$ws.Range("J101").Value = 64489.5
$ws.Range("L101").Value = 64489.5
$ws.Range("N101").Value = -70979.5

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 3023.5789
$ws.Range("I122").Value = 1950.7693
$ws.Range("K122").Value = 5852.3079
$ws.Range("M122").Value = -3402.3079
# Row 132
$ws.Range("H132").Value = 8135.8184
$ws.Range("I132").Value = 7370.5713
$ws.Range("K132").Value = 22111.7139
$ws.Range("L132").Value = 28425
$ws.Range("M132").Value = -19581.7139
# Row 136
$ws.Range("H136").Value = 3335.1482
$ws.Range("I136").Value = 2402
$ws.Range("K136").Value = 7206
$ws.Range("M136").Value = -4656
